$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update AE6: insert new referenced_works URL (W2014085028) into the list ---
$ae6Value = @'
c("https://openalex.org/W149395048", "https://openalex.org/W1524552663", "https://openalex.org/W1912231050", "https://openalex.org/W1965855684", "https://openalex.org/W1967878389", "https://openalex.org/W1971910104", "https://openalex.org/W1975470702", "https://openalex.org/W1991040129", "https://openalex.org/W1996430354", "https://openalex.org/W1999847898", "https://openalex.org/W2014085028", "https://openalex.org/W2045280766", "https://openalex.org/W2075247410", "https://openalex.org/W2082428100", 
"https://openalex.org/W2093127022", "https://openalex.org/W2093191118", "https://openalex.org/W2130276273", "https://openalex.org/W2140363334", "https://openalex.org/W2150260015", "https://openalex.org/W2157751307", "https://openalex.org/W2160172778", "https://openalex.org/W2223777709", "https://openalex.org/W2320659285", "https://openalex.org/W2330823310", "https://openalex.org/W2432120366", "https://openalex.org/W2790072344", "https://openalex.org/W2901173776", "https://openalex.org/W2950183821", 
"https://openalex.org/W2987100386", "https://openalex.org/W3009579572", "https://openalex.org/W3027898185", "https://openalex.org/W3120391982", "https://openalex.org/W3151650706", "https://openalex.org/W3175313553", "https://openalex.org/W3191380389", "https://openalex.org/W3203108148", "https://openalex.org/W4200466116", "https://openalex.org/W4283513283", "https://openalex.org/W4293251519", "https://openalex.org/W4309725708", "https://openalex.org/W4318315097")
'@
$ws.Range("AE6").Formula = "'" + $ae6Value
$ws.Range("AE6").Style = "Normal"
$ws.Rows(6).AutoFit()

# --- Update AF6: replace related_works list entirely ---
$ws.Range("AF6").Formula = "'" + 'c("https://openalex.org/W4206297848", "https://openalex.org/W2234262312", "https://openalex.org/W2039788050", "https://openalex.org/W2562720103", "https://openalex.org/W4245072705", "https://openalex.org/W3128559897", "https://openalex.org/W4247516371", "https://openalex.org/W4248626785", "https://openalex.org/W2162565381", "https://openalex.org/W2064752607")'
$ws.Range("AF6").Style = "Normal"

# --- Add new row 9 (new work record) ---
$ws.Range("A9").Formula = "'" + 'https://openalex.org/W4391533862'
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Formula = "'" + 'STRUCTURE OF REGULATED AGRICULTURAL MARKETS IN INDIA'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Formula = "'" + 'list(au_id = "https://openalex.org/A5017320231", au_display_name = "S. J. Patel", au_orcid = "https://orcid.org/0000-0002-7079-565X", author_position = "first", au_affiliation_raw = "", institution_id = NA, institution_display_name = NA, institution_ror = NA, institution_country_code = NA, institution_type = NA, institution_lineage = NA)'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Formula = "'" + 'This paper explores the structure of regulated agricultural markets in contemporary India. It discusses the evolution of agricultural market regulations from the introduction of APMC Acts after the independence to the recent reforms in agricultural marketing. It critically analyses the arguments made for the deregulation of these markets. Using data from the Situational Assessment of Agricultural Households and Land and Livestock Holdings of Households in Rural India, 2019, the paper also discusses the diverse marketing channels for different crops at the national level.'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'" + '2023-01-01'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Formula = "'" + 'International journal of social science and economic research'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Formula = "'" + 'https://openalex.org/S4210219729'
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Formula = "'" + 'N/A'
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Formula = "'" + '2455-8834'
$ws.Range("I9").Style = "Normal"
$ws.Range("J9").Formula = "'" + 'https://doi.org/10.46609/ijsser.2023.v08i12.018'
$ws.Range("J9").Style = "Normal"
$ws.Range("K9").Formula = "'" + 'N/A'
$ws.Range("K9").Style = "Normal"
$ws.Range("L9").Formula = "'" + 'N/A'
$ws.Range("L9").Style = "Normal"
$ws.Range("M9").Formula = "'" + 'publishedVersion'
$ws.Range("M9").Style = "Normal"
$ws.Range("N9").Formula = "'" + '3963'
$ws.Range("N9").Style = "Normal"
$ws.Range("O9").Formula = "'" + '3973'
$ws.Range("O9").Style = "Normal"
$ws.Range("P9").Formula = "'" + '08'
$ws.Range("P9").Style = "Normal"
$ws.Range("Q9").Formula = "'" + '12'
$ws.Range("Q9").Style = "Normal"
$ws.Range("R9").Formula = "'" + 'TRUE'
$ws.Range("R9").Style = "Normal"
$ws.Range("S9").Formula = "'" + 'TRUE'
$ws.Range("S9").Style = "Normal"
$ws.Range("T9").Formula = "'" + 'bronze'
$ws.Range("T9").Style = "Normal"
$ws.Range("U9").Formula = "'" + 'https://doi.org/10.46609/ijsser.2023.v08i12.018'
$ws.Range("U9").Style = "Normal"
$ws.Range("V9").Formula = "'" + 'FALSE'
$ws.Range("V9").Style = "Normal"
$ws.Range("W9").Formula = "'" + 'en'
$ws.Range("W9").Style = "Normal"
$ws.Range("X9").Formula = "'" + 'N/A'
$ws.Range("X9").Style = "Normal"
$ws.Range("Y9").Formula = "'" + '0'
$ws.Range("Y9").Style = "Normal"
$ws.Range("Z9").Formula = "'" + '2023'
$ws.Range("Z9").Style = "Normal"
$ws.Range("AA9").Formula = "'" + 'https://api.openalex.org/works?filter=cites:W4391533862'
$ws.Range("AA9").Style = "Normal"
$ws.Range("AB9").Formula = "'" + 'c(openalex = "https://openalex.org/W4391533862", doi = "https://doi.org/10.46609/ijsser.2023.v08i12.018")'
$ws.Range("AB9").Style = "Normal"
$ws.Range("AC9").Formula = "'" + 'https://doi.org/10.46609/ijsser.2023.v08i12.018'
$ws.Range("AC9").Style = "Normal"
$ws.Range("AD9").Formula = "'" + 'article'
$ws.Range("AD9").Style = "Normal"
$ws.Range("AE9").Formula = "'" + 'NA'
$ws.Range("AE9").Style = "Normal"
$ws.Range("AF9").Formula = "'" + 'c("https://openalex.org/W3122389410", "https://openalex.org/W2200951064", "https://openalex.org/W2748952813", "https://openalex.org/W2355956201", "https://openalex.org/W2386195957", "https://openalex.org/W2327874825", "https://openalex.org/W2351852648", "https://openalex.org/W2613051533", "https://openalex.org/W2349774843", "https://openalex.org/W2775541961")'
$ws.Range("AF9").Style = "Normal"
$ws.Range("AG9").Formula = "'" + 'FALSE'
$ws.Range("AG9").Style = "Normal"
$ws.Range("AH9").Formula = "'" + 'FALSE'
$ws.Range("AH9").Style = "Normal"
